$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, applied identically to both the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets.
$updates = @{
    3  = 8583
    4  = 1541
    6  = 396
    7  = 278
    9  = 35
    10 = 133
    12 = 466
    13 = 1289
    14 = 327
    15 = 87
    16 = 152
    17 = 106
    18 = 145
    19 = 86
    20 = 130
    21 = 120
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
